$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab26")

# Row 69
$ws.Range("C69").Value = 106503.960907612
$ws.Range("D69").Value = 57840.4577911894
$ws.Range("E69").Value = 48663.5031164226
$ws.Range("F69").Value = 37.1196261276077
$ws.Range("G69").Value = 60.3408010877983
$ws.Range("H69").Value = 25.4696608007859
$ws.Range("N69").Value = 17.9476
$ws.Range("Q69").Value = 280.966
$ws.Range("T69").Value = 5341.549
$ws.Range("U69").Value = 261.384
$ws.Range("Y69").Value = 2539.796
$ws.Range("Z69").Value = 336.41
$ws.Range("AC69").Value = 23822.13
$ws.Range("AH69").Value = 8175.9

# Row 77
$ws.Range("C77").Value = 401951.097636115
$ws.Range("D77").Value = 339210.967991935
$ws.Range("E77").Value = 62740.1296441796
$ws.Range("F77").Value = 99.8029083881313
$ws.Range("G77").Value = 99.9855313910239
$ws.Range("H77").Value = 98.8269791808391
$ws.Range("Q77").Value = 12693.125
$ws.Range("R77").Value = 23960.72
$ws.Range("T77").Value = 160189.779
$ws.Range("Y77").Value = 88653.333
$ws.Range("Z77").Value = 22062.58
$ws.Range("AA77").Value = 79506.09
$ws.Range("AC77").Value = 606959.64
$ws.Range("AH77").Value = 285988.9

# Row 80
$ws.Range("C80").Value = 76718.8307176406
$ws.Range("D80").Value = 63111.2317890846
$ws.Range("E80").Value = 13607.598928556
$ws.Range("F80").Value = 61.151464488815
$ws.Range("G80").Value = 83.2818792434467
$ws.Range("H80").Value = 27.3922650764986
$ws.Range("I80").Value = 10.5569
$ws.Range("L80").Value = ".."
$ws.Range("N80").Value = 27.9395
$ws.Range("Q80").Value = 565.26
$ws.Range("T80").Value = 4628.606
$ws.Range("U80").Value = 66.4
$ws.Range("Y80").Value = 39189.582
$ws.Range("Z80").Value = 840.15
$ws.Range("AC80").Value = 16007.07
$ws.Range("AD80").Value = 210.72
$ws.Range("AH80").Value = 123102.06

# Row 82
$ws.Range("C82").Value = 685432.726857059
$ws.Range("D82").Value = 423225.145328048
$ws.Range("E82").Value = 262207.581529011
$ws.Range("F82").Value = 55.5728347480012
$ws.Range("G82").Value = 82.0736341456908
$ws.Range("H82").Value = 36.5328936300961
$ws.Range("I82").Value = 532.4957
$ws.Range("L82").Value = 850.7815
$ws.Range("N82").Value = 203.0789
$ws.Range("Q82").Value = 10000.259
$ws.Range("T82").Value = 28962.715
$ws.Range("U82").Value = 1585.812
$ws.Range("Y82").Value = 138253.416
$ws.Range("Z82").Value = 15924.24
$ws.Range("AC82").Value = 119722.65
$ws.Range("AD82").Value = 2751.73
$ws.Range("AH82").Value = 497337.37

# Row 84
$ws.Range("C84").Value = 208227.698432706
$ws.Range("D84").Value = 126092.872874768
$ws.Range("E84").Value = 82134.8255579385
$ws.Range("F84").Value = 35.8216428126844
$ws.Range("G84").Value = 67.318390259981
$ws.Range("H84").Value = 20.8473555628059
$ws.Range("I84").Value = 268.0345
$ws.Range("L84").Value = 365.6051
$ws.Range("N84").Value = 144.4334
$ws.Range("Q84").Value = 782.36
$ws.Range("T84").Value = 15551.071
$ws.Range("U84").Value = 738.729
$ws.Range("Y84").Value = 7005.574
$ws.Range("Z84").Value = 927.22
$ws.Range("AC84").Value = 65549.52
$ws.Range("AD84").Value = 618.51
$ws.Range("AH84").Value = 16635.78

# Row 86
$ws.Range("C86").Value = 491868.002775684
$ws.Range("D86").Value = 314242.363964352
$ws.Range("E86").Value = 177625.638811332
$ws.Range("F86").Value = 70.0887513087416
$ws.Range("G86").Value = 89.2123643177245
$ws.Range("H86").Value = 50.8172632495077
$ws.Range("I86").Value = 274.9588
$ws.Range("L86").Value = 485.1764
$ws.Range("N86").Value = 79.8592
$ws.Range("Q86").Value = 3548.891
$ws.Range("T86").Value = 16423.433
$ws.Range("U86").Value = 556.698
$ws.Range("Y86").Value = 114167.773
$ws.Range("Z86").Value = 7621.28
$ws.Range("AC86").Value = 66485.29
$ws.Range("AD86").Value = 1502.17
$ws.Range("AH86").Value = 375573.38

# Row 87
$ws.Range("C87").Value = 2528102.30508261
$ws.Range("D87").Value = 1064290.4216536
$ws.Range("E87").Value = 1463811.88342901
$ws.Range("F87").Value = 95.4271185701114
$ws.Range("G87").Value = 99.4262531781985
$ws.Range("H87").Value = 92.7157209308635
$ws.Range("I87").Value = 975.1429
$ws.Range("N87").Value = 585.0373
$ws.Range("Q87").Value = 68398.552
$ws.Range("R87").Value = 43457.749
$ws.Range("T87").Value = 133793.92
$ws.Range("U87").Value = 15209.846
$ws.Range("Y87").Value = 599989.286
$ws.Range("Z87").Value = 85528.71
$ws.Range("AA87").Value = 75084.54
$ws.Range("AC87").Value = 449724.37
$ws.Range("AD87").Value = 43513.11
$ws.Range("AH87").Value = 2388852.44

# Row 89
$ws.Range("C89").Value = 2433266.98955843
$ws.Range("D89").Value = 1643773.69934272
$ws.Range("E89").Value = 789305.564405691
$ws.Range("F89").Value = 99.9254460378116
$ws.Range("G89").Value = 99.9637712125535
$ws.Range("H89").Value = 99.8220035418072
$ws.Range("I89").Value = 986.4707
$ws.Range("N89").Value = 2005.9269
$ws.Range("Q89").Value = 286501.702
$ws.Range("R89").Value = 323801.895
$ws.Range("S89").Value = 2956.948
$ws.Range("T89").Value = 620782.724
$ws.Range("U89").Value = 51421.225
$ws.Range("W89").Value = 86995.5
$ws.Range("X89").Value = 35348
$ws.Range("Y89").Value = 1764453.055
$ws.Range("Z89").Value = 312753.7
$ws.Range("AA89").Value = 594080.8
$ws.Range("AC89").Value = 2316769.53
$ws.Range("AD89").Value = 210946.88
$ws.Range("AF89").Value = 636306.43
$ws.Range("AG89").Value = 37746.7
$ws.Range("AH89").Value = 7289185.21

# Row 90
$ws.Range("C90").Value = 1204475.52304965
$ws.Range("D90").Value = 977707.540061529
$ws.Range("E90").Value = 226767.982988118
$ws.Range("F90").Value = 99.9882522041636
$ws.Range("G90").Value = 99.999054632731
$ws.Range("H90").Value = 99.941704369013
$ws.Range("N90").Value = 506.1557
$ws.Range("Q90").Value = 344046.364
$ws.Range("R90").Value = 356677.138
$ws.Range("S90").Value = 5718.973
$ws.Range("T90").Value = 401769.61
$ws.Range("U90").Value = 64158.524
$ws.Range("W90").Value = 286704.09
$ws.Range("X90").Value = 74805.74
$ws.Range("Y90").Value = 1722608.394
$ws.Range("Z90").Value = 420830.05
$ws.Range("AA90").Value = 905541.11
$ws.Range("AC90").Value = 1382177.89
$ws.Range("AD90").Value = 323597.34
$ws.Range("AF90").Value = 1861470.66
$ws.Range("AG90").Value = 66589.73
$ws.Range("AH90").Value = 6064958.93

# Row 91
$ws.Range("C91").Value = 271821.589462866

# Row 94
$ws.Range("C94").Value = 50774.9713313152
$ws.Range("D94").Value = 35298.1996258911
$ws.Range("E94").Value = 15289.0458954061
$ws.Range("F94").Value = 83.1113640713163
$ws.Range("G94").Value = 95.5812003008027
$ws.Range("H94").Value = 63.2755899438685
$ws.Range("Q94").Value = 1217.705
$ws.Range("R94").Value = 505.702
$ws.Range("Y94").Value = 28871.775
$ws.Range("Z94").Value = 1300.95
$ws.Range("AA94").Value = 1429.05
$ws.Range("AH94").Value = 111299.07

# Row 97
$ws.Range("C97").Value = 395922.054453626
$ws.Range("D97").Value = 263086.295496786
$ws.Range("E97").Value = 132835.75895684
$ws.Range("F97").Value = 49.483357342189
$ws.Range("G97").Value = 81.8178606724979
$ws.Range("H97").Value = 27.7573674354215
$ws.Range("I97").Value = 298.4385
$ws.Range("M97").Value = 8.9954
$ws.Range("N97").Value = 197.2817
$ws.Range("P97").Value = 839.8455
$ws.Range("Q97").Value = 794.972
$ws.Range("R97").Value = 698.297
$ws.Range("T97").Value = 20977.41
$ws.Range("U97").Value = 942.658
$ws.Range("Y97").Value = 34088.778
$ws.Range("Z97").Value = 978.45
$ws.Range("AA97").Value = 792.19
$ws.Range("AC97").Value = 82222
$ws.Range("AD97").Value = 1447.25
$ws.Range("AH97").Value = 94837.87

# Row 98
$ws.Range("C98").Value = 576436.726021037
$ws.Range("D98").Value = 293498.145096111
$ws.Range("E98").Value = 282938.580924926
$ws.Range("F98").Value = 85.3752532168483
$ws.Range("G98").Value = 95.9714766719299
$ws.Range("H98").Value = 76.6019716273489
$ws.Range("I98").Value = 239.8361
$ws.Range("N98").Value = 145.6779
$ws.Range("P98").Value = 775.6098
$ws.Range("Q98").Value = 2650.824
$ws.Range("R98").Value = 2153.271
$ws.Range("T98").Value = 58261.931
$ws.Range("U98").Value = 2150.704
$ws.Range("X98").Value = 1040
$ws.Range("Y98").Value = 156040.319
$ws.Range("Z98").Value = 3768.16
$ws.Range("AA98").Value = 4724.08
$ws.Range("AC98").Value = 202220.62
$ws.Range("AD98").Value = 8502.26
$ws.Range("AH98").Value = 568966.93

